# Update MVX codeset workbook:
#  - refresh three "last updated date" values to 44866 (2022-11-01)
#  - append a new manufacturer row (CAN / CanSino Biologics, Inc)
#  - grow the query table / defined name to cover the new row
#  - update the table / query names (drop stale "_1" suffix)
#  - move the active selection as recorded by the saved session

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Correct existing "last updated date" values (col E) ----
$ws.Range("E81").Value = 44866
$ws.Range("E82").Value = 44866
$ws.Range("E83").Value = 44866

# ---- 2. Grow the table by one row and populate it ----
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Range("A90").Value = "CAN"
$ws.Range("B90").Value = "CanSino Biologics, Inc"
$ws.Range("C90").Value = "Non-US COVID-19 vaccine manufacturer,  WHO Authorized (CONVIDECIA) on 5/19/2022"
$ws.Range("D90").Value = "Active"
$ws.Range("E90").Value = 44866
$ws.Range("F90").Value = 99

# ---- 3. Rename table / query table, dropping the stale "_1" suffix ----
$lo.Name = "Table_DSDV_INFC_1601_qsrv1_NIP_INSIDENIP_tblMVXCodes"
$qt = $lo.QueryTable
$qt.Name = "DSDV-INFC-1601_qsrv1 NIP_INSIDENIP tblMVXCodes"

# ---- 4. Update the hidden defined name that tracks the query range ----
$nm = $wb.Names.Item(1)
$nm.Name = "DSDV_INFC_1601_qsrv1_NIP_INSIDENIP_tblMVXCodes"
$nm.RefersTo = "='WEB mvx'!`$A`$1:`$F`$90"

# ---- 5. Restore the recorded selection ----
$ws.Range("B13").Select() | Out-Null
